# Proteomics-Template.xlsx: update proteomics handlers to the new format
# template from Rolf (CCS-7, SP-30).
#
# Content changes applied:
#   1. "Value Unit" description (openbis-metadata!C6) gains a new allowed
#      unit, "fmol/ug protein digest".
#   2. The active selection on the metadata sheet moves from C8 to C7.
#   3. Column C is widened (best-fit) to show the longer description text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# 1. Extend the list of valid "Value Unit" values shown in the Description
#    column next to "Value Unit" (row 6).
$ws.Range("C6").Value = "One of mM, uM, Percent, RatioT1, RatioCs, or AU, Dimensionless, fmol/ug protein digest"

# 2. Widen column C so the longer description text fits/best-fits.
$ws.Columns.Item(3).ColumnWidth = 89.28571428571429

# 3. Move the active cell/selection to C7 (was C8).
$ws.Range("C7").Select()
